$d = $word.ActiveDocument

# The "Metodos" table cell contains a duplicated block:
#   "Instalação de Dashboards" / "" / "Instalar dispositivos..." /
#   "Treinar os moradores..." / ""
# repeated twice in a row, immediately followed by "Gamificação".
# The edit removes the first (duplicate) copy of that 5-paragraph block,
# leaving a single copy before "Gamificação".

function Get-ParaText($idx) {
    return $d.Paragraphs($idx).Range.Text.TrimEnd([char]13)
}

$headingText = "Instalação de Dashboards"

# Locate the paragraph index of the first occurrence of the heading.
$firstIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ((Get-ParaText $i) -eq $headingText) {
        $firstIndex = $i
        break
    }
}

if ($firstIndex -eq -1) {
    throw "Could not locate 'Instalação de Dashboards' paragraph."
}

# Confirm the expected duplicated-block shape before touching anything:
#   firstIndex+0 : Instalação de Dashboards
#   firstIndex+1 : (empty)
#   firstIndex+2 : Instalar dispositivos...
#   firstIndex+3 : Treinar os moradores...
#   firstIndex+4 : (empty)
#   firstIndex+5 : Instalação de Dashboards   <- second (kept) copy starts here
$expectedSecond = $firstIndex + 5
if ((Get-ParaText $expectedSecond) -ne $headingText) {
    throw "Unexpected document structure; duplicated heading not found where expected."
}

# Delete the first copy's 5 paragraphs (indices firstIndex .. firstIndex+4),
# working from the last one back to the first so earlier indices stay valid
# as each paragraph is removed.
for ($i = $firstIndex + 4; $i -ge $firstIndex; $i--) {
    $d.Paragraphs($i).Range.Delete()
}
